$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3371.75
$ws.Range("I74").Value = 2992.5
$ws.Range("J74").Value = 3751
$ws.Range("K74").Value = 2992.5
$ws.Range("L74").Value = 3751
$ws.Range("M74").Value = -2056.5
$ws.Range("N74").Value = -5623

$ws.Range("H77").Value = 3371.75
$ws.Range("I77").Value = 2992.5
$ws.Range("J77").Value = 3751
$ws.Range("K77").Value = 14962.5
$ws.Range("L77").Value = 18755
$ws.Range("M77").Value = -10282.5
$ws.Range("N77").Value = -28115

$ws.Range("H98").Value = 3247.0527
$ws.Range("I98").Value = 1846.2667
$ws.Range("J98").Value = 8500
$ws.Range("K98").Value = 1846.2667
$ws.Range("L98").Value = 8500
$ws.Range("M98").Value = -348.2666999999999
$ws.Range("N98").Value = -11496

$ws.Range("H113").Value = 4420.343
$ws.Range("I113").Value = 3630.7693
$ws.Range("J113").Value = 4886.909
$ws.Range("K113").Value = 3630.7693
$ws.Range("L113").Value = 4886.909
$ws.Range("M113").Value = -376.7692999999999
$ws.Range("N113").Value = -11394.909

$ws.Range("H116").Value = 3083.8333
$ws.Range("I116").Value = 2997
$ws.Range("K116").Value = 2997
$ws.Range("M116").Value = 445

$ws.Range("H122").Value = 3247.0527
$ws.Range("I122").Value = 1846.2667
$ws.Range("J122").Value = 8500
$ws.Range("K122").Value = 5538.800099999999
$ws.Range("L122").Value = 25500
$ws.Range("M122").Value = -3088.800099999999
$ws.Range("N122").Value = -30400

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8803.672
$ws.Range("I32").Value = 3896.3333
$ws.Range("K32").Value = 3896.3333
$ws.Range("M32").Value = -3609.3333

$ws.Range("H63").Value = 3916.3635
$ws.Range("J63").Value = 4714.2856
$ws.Range("L63").Value = 4714.2856
$ws.Range("N63").Value = -6086.2856

$ws.Range("H66").Value = 3916.3635
$ws.Range("J66").Value = 4714.2856
$ws.Range("L66").Value = 23571.428
$ws.Range("N66").Value = -30435.428

$ws.Range("H88").Value = 1562.375
$ws.Range("I88").Value = 1683.3334
$ws.Range("J88").Value = 1199.5
$ws.Range("K88").Value = 1683.3334
$ws.Range("L88").Value = 1199.5
$ws.Range("M88").Value = -1277.3334
$ws.Range("N88").Value = -2011.5

$ws.Range("H91").Value = 1562.375
$ws.Range("I91").Value = 1683.3334
$ws.Range("J91").Value = 1199.5
$ws.Range("K91").Value = 1683.3334
$ws.Range("L91").Value = 1199.5
$ws.Range("M91").Value = -279.3334
$ws.Range("N91").Value = -4007.5

$ws.Range("H137").Value = 38600
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 38600
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 38600
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -48800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1715.4474
$ws.Range("I86").Value = 1674.375
$ws.Range("J86").Value = 1745.3182
$ws.Range("K86").Value = 1674.375
$ws.Range("L86").Value = 1745.3182
$ws.Range("M86").Value = -551.375
$ws.Range("N86").Value = -3991.3182

$ws.Range("H89").Value = 1715.4474
$ws.Range("I89").Value = 1674.375
$ws.Range("J89").Value = 1745.3182
$ws.Range("K89").Value = 8371.875
$ws.Range("L89").Value = 8726.591
$ws.Range("M89").Value = -2755.875
$ws.Range("N89").Value = -19958.591

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 881.0833
$ws.Range("I122").Value = 902.55554
$ws.Range("J122").Value = 816.6667
$ws.Range("K122").Value = 2707.66662
$ws.Range("L122").Value = 2450.0001
$ws.Range("M122").Value = -257.66662
$ws.Range("N122").Value = -7350.0001

$ws.Range("H141").Value = 67804.28999999999
$ws.Range("J141").Value = 67804.28999999999
$ws.Range("L141").Value = 67804.28999999999
$ws.Range("N141").Value = -78164.28999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 2620
$ws.Range("I19").Value = 100
$ws.Range("J19").Value = 4300
$ws.Range("K19").Value = 300
$ws.Range("L19").Value = 12900
$ws.Range("M19").Value = -126
$ws.Range("N19").Value = -13248

$ws.Range("H131").Value = 912.38
$ws.Range("I131").Value = 299.6
$ws.Range("J131").Value = 944.6316
$ws.Range("K131").Value = 898.8000000000001
$ws.Range("L131").Value = 2833.8948
$ws.Range("M131").Value = 4141.2
$ws.Range("N131").Value = -12913.8948

$ws.Range("H140").Value = 2305
$ws.Range("I140").Value = 735.4375
$ws.Range("K140").Value = 2206.3125
$ws.Range("M140").Value = 2973.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2313
$ws.Range("I80").Value = 2400
$ws.Range("J80").Value = 2052
$ws.Range("K80").Value = 2400
$ws.Range("L80").Value = 2052
$ws.Range("M80").Value = -1402
$ws.Range("N80").Value = -4048

$ws.Range("H83").Value = 2313
$ws.Range("I83").Value = 2400
$ws.Range("J83").Value = 2052
$ws.Range("K83").Value = 12000
$ws.Range("L83").Value = 10260
$ws.Range("M83").Value = -7008
$ws.Range("N83").Value = -20244

$ws.Range("H102").Value = 753.86487
$ws.Range("I102").Value = 747.6
$ws.Range("J102").Value = 863.5
$ws.Range("K102").Value = 747.6
$ws.Range("L102").Value = 863.5
$ws.Range("M102").Value = 874.4
$ws.Range("N102").Value = -4107.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1507.258
$ws.Range("I40").Value = 1497.6923
$ws.Range("J40").Value = 1557
$ws.Range("K40").Value = 1497.6923
$ws.Range("L40").Value = 1557
$ws.Range("M40").Value = -1361.6923
$ws.Range("N40").Value = -1829

$ws.Range("H82").Value = 1225.1765
$ws.Range("I82").Value = 992
$ws.Range("J82").Value = 1388.4
$ws.Range("K82").Value = 992
$ws.Range("L82").Value = 1388.4
$ws.Range("M82").Value = -631
$ws.Range("N82").Value = -2110.4

$ws.Range("H85").Value = 1225.1765
$ws.Range("I85").Value = 992
$ws.Range("J85").Value = 1388.4
$ws.Range("K85").Value = 992
$ws.Range("L85").Value = 1388.4
$ws.Range("M85").Value = 256
$ws.Range("N85").Value = -3884.4

$ws.Range("H122").Value = 2911.2173
$ws.Range("I122").Value = 2144.2222
$ws.Range("J122").Value = 3404.2856
$ws.Range("K122").Value = 6432.6666
$ws.Range("L122").Value = 10212.8568
$ws.Range("M122").Value = -3982.6666
$ws.Range("N122").Value = -15112.8568

$ws.Range("H132").Value = 3933.818
$ws.Range("I132").Value = 3276.1333
$ws.Range("J132").Value = 5343.143
$ws.Range("K132").Value = 9828.3999
$ws.Range("L132").Value = 16029.429
$ws.Range("M132").Value = -7298.3999
$ws.Range("N132").Value = -21089.429

$ws.Range("H136").Value = 6632.8945
$ws.Range("I136").Value = 1401.3334
$ws.Range("J136").Value = 26251.25
$ws.Range("K136").Value = 4204.0002
$ws.Range("L136").Value = 78753.75
$ws.Range("M136").Value = -1654.0002
$ws.Range("N136").Value = -83853.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1364.9
$ws.Range("I81").Value = 1249.75
$ws.Range("J81").Value = 1441.6666
$ws.Range("K81").Value = 2499.5
$ws.Range("L81").Value = 2883.3332
$ws.Range("M81").Value = -1438.5
$ws.Range("N81").Value = -5005.3332

$ws.Range("H84").Value = 1364.9
$ws.Range("I84").Value = 1249.75
$ws.Range("J84").Value = 1441.6666
$ws.Range("K84").Value = 12497.5
$ws.Range("L84").Value = 14416.666
$ws.Range("M84").Value = -7193.5
$ws.Range("N84").Value = -25024.666

$ws.Range("H88").Value = 27180
$ws.Range("J88").Value = 27180
$ws.Range("L88").Value = 27180
$ws.Range("N88").Value = -27992

$ws.Range("H91").Value = 27180
$ws.Range("J91").Value = 27180
$ws.Range("L91").Value = 27180
$ws.Range("N91").Value = -29988

$ws.Range("H122").Value = 80990.53
$ws.Range("I122").Value = 1075
$ws.Range("J122").Value = 134267.56
$ws.Range("K122").Value = 3225
$ws.Range("L122").Value = 402802.68
$ws.Range("M122").Value = -775
$ws.Range("N122").Value = -407702.68
